$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map country name (column A) -> region
$regionMap = @{
    "Australia"       = "Oceania"
    "Austria"         = "Europe"
    "Belgium"         = "Europe"
    "Canada"          = "North America"
    "Chile"           = "Latin America"
    "Czech Republic"  = "Europe"
    "Germany"         = "Europe"
    "Denmark"         = "Europe"
    "United Kingdom"  = "Europe"
    "Spain"           = "Europe"
    "Finland"         = "Europe"
    "France"          = "Europe"
    "Greece"          = "Europe"
    "Hungary"         = "Europe"
    "Ireland"         = "Europe"
    "Israel"          = "Europe"
    "Italy"           = "Europe"
    "Japan"           = "Asia"
    "Korea"           = "Asia"
    "Luxembourg"      = "Europe"
    "Latvia"          = "Europe"
    "Mexico"          = "Latin America"
    "Netherlands"     = "Europe"
    "Norway"          = "Europe"
    "New Zealand"     = "Oceania"
    "Poland"          = "Europe"
    "Portugal"        = "Europe"
    "Slovak Republic" = "Europe"
    "Slovenia"        = "Europe"
    "Sweden"          = "Europe"
    "Turkey"          = "Europe"
    "United States"   = "North America"
    "Colombia"        = "Latin America"
    "Costa Rica"      = "Latin America"
    "Lithuania"       = "Europe"
}

# Header
$ws.Cells.Item(1, 12).Value = "region"

# Fill region values for each data row (2..36) based on country in column A
for ($r = 2; $r -le 36; $r++) {
    $country = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 12).Value = $regionMap[$country]
}

# Set column L width to match bestFit/customWidth sizing used in the target file
# (13.6 in "character width" units renders as width="14.5" in the underlying OOXML)
$ws.Columns.Item(12).ColumnWidth = 13.6

# Update view to reflect scrolled position / selection seen in target file
$excel.Goto($ws.Range("D1"), $true)
$ws.Range("L10").Select()
